{"js": "// Update the header date and every three-digit x one-digit multiplication\n// problem in the table to the new set of values.\nconst replacements = [\n  [\"2025-12-02 Tuesday\", \"2025-12-03 Wednesday\"],\n  [\"590\u00d79=5310\", \"812\u00d75=4060\"],\n  [\"271\u00d72=542\", \"454\u00d76=2724\"],\n  [\"350\u00d73=1050\", \"344\u00d72=688\"],\n  [\"842\u00d74=3368\", \"126\u00d75=630\"],\n  [\"603\u00d75=3015\", \"285\u00d73=855\"],\n  [\"163\u00d72=326\", \"414\u00d79=3726\"],\n  [\"251\u00d76=1506\", \"377\u00d77=2639\"],\n  [\"204\u00d79=1836\", \"845\u00d79=7605\"],\n  [\"526\u00d73=1578\", \"823\u00d75=4115\"],\n  [\"178\u00d77=1246\", \"638\u00d78=5104\"],\n  [\"118\u00d79=1062\", \"544\u00d73=1632\"],\n  [\"923\u00d75=4615\", \"334\u00d77=2338\"],\n  [\"608\u00d77=4256\", \"788\u00d74=3152\"],\n  [\"156\u00d75=780\", \"585\u00d78=4680\"],\n  [\"131\u00d77=917\", \"908\u00d76=5448\"],\n  [\"882\u00d78=7056\", \"809\u00d72=1618\"],\n  [\"186\u00d76=1116\", \"334\u00d73=1002\"],\n  [\"155\u00d75=775\", \"348\u00d79=3132\"],\n  [\"795\u00d74=3180\", \"683\u00d72=1366\"],\n  [\"571\u00d76=3426\", \"208\u00d73=624\"],\n  [\"314\u00d72=628\", \"898\u00d76=5388\"],\n  [\"546\u00d73=1638\", \"353\u00d74=1412\"],\n  [\"979\u00d72=1958\", \"945\u00d78=7560\"],\n  [\"749\u00d73=2247\", \"126\u00d76=756\"],\n  [\"760\u00d78=6080\", \"197\u00d79=1773\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the header date and every three-digit x one-digit multiplication\n# problem in the table to the new set of values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-12-02 Tuesday\", \"2025-12-03 Wednesday\"),\n  @(\"590\u00d79=5310\", \"812\u00d75=4060\"),\n  @(\"271\u00d72=542\", \"454\u00d76=2724\"),\n  @(\"350\u00d73=1050\", \"344\u00d72=688\"),\n  @(\"842\u00d74=3368\", \"126\u00d75=630\"),\n  @(\"603\u00d75=3015\", \"285\u00d73=855\"),\n  @(\"163\u00d72=326\", \"414\u00d79=3726\"),\n  @(\"251\u00d76=1506\", \"377\u00d77=2639\"),\n  @(\"204\u00d79=1836\", \"845\u00d79=7605\"),\n  @(\"526\u00d73=1578\", \"823\u00d75=4115\"),\n  @(\"178\u00d77=1246\", \"638\u00d78=5104\"),\n  @(\"118\u00d79=1062\", \"544\u00d73=1632\"),\n  @(\"923\u00d75=4615\", \"334\u00d77=2338\"),\n  @(\"608\u00d77=4256\", \"788\u00d74=3152\"),\n  @(\"156\u00d75=780\", \"585\u00d78=4680\"),\n  @(\"131\u00d77=917\", \"908\u00d76=5448\"),\n  @(\"882\u00d78=7056\", \"809\u00d72=1618\"),\n  @(\"186\u00d76=1116\", \"334\u00d73=1002\"),\n  @(\"155\u00d75=775\", \"348\u00d79=3132\"),\n  @(\"795\u00d74=3180\", \"683\u00d72=1366\"),\n  @(\"571\u00d76=3426\", \"208\u00d73=624\"),\n  @(\"314\u00d72=628\", \"898\u00d76=5388\"),\n  @(\"546\u00d73=1638\", \"353\u00d74=1412\"),\n  @(\"979\u00d72=1958\", \"945\u00d78=7560\"),\n  @(\"749\u00d73=2247\", \"126\u00d76=756\"),\n  @(\"760\u00d78=6080\", \"197\u00d79=1773\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute([ref]$oldText, $false, $true, $false, $false, $false, $true, 0, $false, [ref]$newText, 2) | Out-Null\n}\n"}
